$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "Move to location (11, 8) and remove the toolkit."
$ws.Cells.Item(2, 2).Value = 69.86899699999999
$ws.Cells.Item(2, 3).Value = 10337
$ws.Cells.Item(2, 4).Value = "'0.03084"
$ws.Cells.Item(2, 5).Value = "1436cb36-0897-4d6d-abcf-c5a0e375c0d7"

# Row 3
$ws.Cells.Item(3, 1).Value = "Move to location (7, 5) and remove the liquid spill."
$ws.Cells.Item(3, 2).Value = 69.502499
$ws.Cells.Item(3, 3).Value = 10057
$ws.Cells.Item(3, 4).Value = "'0.0303"
$ws.Cells.Item(3, 5).Value = "28e7e071-82d1-44fc-86c7-f8fab9f30c03"

# Row 4
$ws.Cells.Item(4, 1).Value = "Move to location (8, 6) and remove the large debris."
$ws.Cells.Item(4, 2).Value = 74.52652500000001
$ws.Cells.Item(4, 3).Value = 10035
$ws.Cells.Item(4, 4).Value = "'0.03024"
$ws.Cells.Item(4, 5).Value = "c99f35ad-283b-49ca-8f96-a780164c8567"

# Row 5
$ws.Cells.Item(5, 1).Value = "Move to location (2, 4) and remove the dust."
$ws.Cells.Item(5, 2).Value = 72.06699999999999
$ws.Cells.Item(5, 3).Value = 10477
$ws.Cells.Item(5, 4).Value = "'0.03204"
$ws.Cells.Item(5, 5).Value = "1881dc6a-7a58-4f7a-bd0b-5779cd3803ed"

# Row 6
$ws.Cells.Item(6, 1).Value = "Move to location (5, 2) and remove the grass."
$ws.Cells.Item(6, 2).Value = 72.90800900000001
$ws.Cells.Item(6, 3).Value = 9936
$ws.Cells.Item(6, 4).Value = "'0.03051"
$ws.Cells.Item(6, 5).Value = "cafd0021-f8e5-4905-8ddf-8b61327535d5"

# Row 7
$ws.Cells.Item(7, 1).Value = "Move to location (6, 7) and remove the small debris."
$ws.Cells.Item(7, 2).Value = 69.397496
$ws.Cells.Item(7, 3).Value = 9694
$ws.Cells.Item(7, 4).Value = "'0.02667"
$ws.Cells.Item(7, 5).Value = "ccc1397c-31f8-482e-9304-768d60bf4f8a"

# Row 8
$ws.Cells.Item(8, 1).Value = "Move to location (3, 6) and remove the vehicle."
$ws.Cells.Item(8, 2).Value = 72.98699999999999
$ws.Cells.Item(8, 3).Value = 10359
$ws.Cells.Item(8, 4).Value = "'0.0312"
$ws.Cells.Item(8, 5).Value = "ab1cc93b-28a9-4000-bb17-fbda0e764a1c"

# Row 9
$ws.Cells.Item(9, 1).Value = "Move to location (6, 6) and remove the construction materials."
$ws.Cells.Item(9, 2).Value = 69.1215
$ws.Cells.Item(9, 3).Value = 10182
$ws.Cells.Item(9, 4).Value = "'0.02886"
$ws.Cells.Item(9, 5).Value = "5b2996bd-f14a-42c8-9a43-516aaa54c78e"

# Row 10
$ws.Cells.Item(10, 1).Value = "Move to location (3, 9) and remove the tree branches."
$ws.Cells.Item(10, 2).Value = 66.399998
$ws.Cells.Item(10, 3).Value = 10299
$ws.Cells.Item(10, 4).Value = "'0.03312"
$ws.Cells.Item(10, 5).Value = "2e7c5e9d-3e12-4334-911e-406ed60c30c9"

# Row 11
$ws.Cells.Item(11, 1).Value = "Move to location (6, 6) and remove the screws."
$ws.Cells.Item(11, 2).Value = 64.45051599999999
$ws.Cells.Item(11, 3).Value = 10228
$ws.Cells.Item(11, 4).Value = "'0.02985"
$ws.Cells.Item(11, 5).Value = "74fcfe22-19f2-4414-9fbf-fedc381c215f"
